$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 2.38
$ws.Range("W2").Value = 12
$ws.Range("AB2").Value = 23
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 41
$ws.Range("AQ2").Value = 41
$ws.Range("AR2").Value = 51
$ws.Range("AW2").Value = 351
# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("I3").Value = 2.8
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 19
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 17
$ws.Range("AW3").Value = 251
$ws.Range("BB3").Value = 51
# Row 4
$ws.Range("G4").Value = 3.75
$ws.Range("I4").Value = 1.95
$ws.Range("AO4").Value = 21
# Row 5
$ws.Range("Q5").Value = 2.02
$ws.Range("R5").Value = 1.88
# Row 6
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.75
# Row 8
$ws.Range("H8").Value = 4.75
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 1.91
$ws.Range("L8").Value = 9.5
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 8.5
$ws.Range("AD8").Value = 9.5
$ws.Range("AF8").Value = 126
$ws.Range("AJ8").Value = 29
$ws.Range("AK8").Value = 126
$ws.Range("AL8").Value = 81
$ws.Range("AN8").Value = 3
$ws.Range("AV8").Value = 101
$ws.Range("BA8").Value = 301
$ws.Range("BB8").Value = 351
# Row 9
$ws.Range("G9").Value = 1.85
$ws.Range("I9").Value = 4.5
$ws.Range("J9").Value = 2.63
$ws.Range("L9").Value = 5
$ws.Range("Q9").Value = 2.4
$ws.Range("R9").Value = 1.53
$ws.Range("AA9").Value = 19
$ws.Range("AH9").Value = 9.5
$ws.Range("AJ9").Value = 15
$ws.Range("AO9").Value = 11
$ws.Range("AQ9").Value = 41
$ws.Range("AX9").Value = 6
$ws.Range("AY9").Value = 26
# Row 10
$ws.Range("G10").Value = 1.75
$ws.Range("I10").Value = 5.5
$ws.Range("J10").Value = 2.5
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("Q10").Value = 2.5
$ws.Range("R10").Value = 1.5
$ws.Range("U10").Value = 2.25
$ws.Range("V10").Value = 1.57
$ws.Range("X10").Value = 7
$ws.Range("AA10").Value = 17
$ws.Range("AC10").Value = 6.5
$ws.Range("AE10").Value = 21
$ws.Range("AJ10").Value = 19
$ws.Range("AK10").Value = 51
$ws.Range("AM10").Value = 51
$ws.Range("AN10").Value = 3.5
$ws.Range("AO10").Value = 10
$ws.Range("AS10").Value = 251
$ws.Range("AX10").Value = 6.5
$ws.Range("BA10").Value = 126
# Row 12
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 1.93
$ws.Range("R12").Value = 1.93
# Row 13
$ws.Range("Q13").Value = 1.8
$ws.Range("R13").Value = 2
$ws.Range("Z13").Value = 8.5
$ws.Range("AF13").Value = 67
$ws.Range("AG13").Value = 451
$ws.Range("AJ13").Value = 23
$ws.Range("AL13").Value = 51
$ws.Range("AM13").Value = 51
$ws.Range("AU13").Value = 9.5
$ws.Range("AX13").Value = 9
$ws.Range("BC13").Value = 401
# Row 14
$ws.Range("N14").Value = 13
# Row 17
$ws.Range("I17").Value = 11
$ws.Range("X17").Value = 6.5
$ws.Range("AB17").Value = 29
$ws.Range("AH17").Value = 26
$ws.Range("AI17").Value = 51
$ws.Range("AK17").Value = 151
$ws.Range("AL17").Value = 81
$ws.Range("AM17").Value = 67
$ws.Range("AN17").Value = 3.2
$ws.Range("AX17").Value = 12
$ws.Range("BA17").Value = 301
$ws.Range("BB17").Value = 251
# Row 20
$ws.Range("G20").Value = 2.62
$ws.Range("H20").Value = 2.95
$ws.Range("I20").Value = 2.8
$ws.Range("J20").Value = 3.25
$ws.Range("K20").Value = 1.98
$ws.Range("L20").Value = 3.4
$ws.Range("M20").Value = 1.1
$ws.Range("N20").Value = 5.8
$ws.Range("O20").Value = 1.45
$ws.Range("P20").Value = 2.55
$ws.Range("Q20").Value = 2.32
$ws.Range("R20").Value = 1.53
$ws.Range("S20").Value = 1.5
$ws.Range("T20").Value = 2.42
$ws.Range("U20").Value = 1.95
$ws.Range("V20").Value = 1.75
$ws.Range("W20").Value = 6.9
$ws.Range("X20").Value = 12
$ws.Range("Y20").Value = 10
$ws.Range("Z20").Value = 30
$ws.Range("AA20").Value = 25
$ws.Range("AB20").Value = 40
$ws.Range("AC20").Value = 5.8
$ws.Range("AD20").Value = 5.7
$ws.Range("AE20").Value = 16
$ws.Range("AF20").Value = 90
$ws.Range("AG20").Value = 900
$ws.Range("AH20").Value = 7.1
$ws.Range("AI20").Value = 13
$ws.Range("AJ20").Value = 10.5
$ws.Range("AK20").Value = 32
$ws.Range("AL20").Value = 28
$ws.Range("AM20").Value = 40
$ws.Range("AO20").Value = 14.5
$ws.Range("AP20").Value = 24
$ws.Range("AQ20").Value = 65
$ws.Range("AR20").Value = 110
$ws.Range("AS20").Value = 350
$ws.Range("AT20").Value = 2.42
$ws.Range("AU20").Value = 7.3
$ws.Range("AV20").Value = 75
$ws.Range("AX20").Value = 4.6
$ws.Range("AZ20").Value = 25
$ws.Range("BA20").Value = 75
$ws.Range("BB20").Value = 120
$ws.Range("BC20").Value = 350
# Row 25
$ws.Range("G25").Value = 1.25
$ws.Range("H25").Value = 5.1
$ws.Range("J25").Value = 1.65
$ws.Range("K25").Value = 2.5
$ws.Range("N25").Value = 16
$ws.Range("O25").Value = 1.19
$ws.Range("P25").Value = 3.75
$ws.Range("Q25").Value = 1.57
$ws.Range("R25").Value = 2.12
$ws.Range("S25").Value = 1.29
$ws.Range("T25").Value = 3.4
$ws.Range("W25").Value = 7
$ws.Range("X25").Value = 5.9
$ws.Range("Y25").Value = 9
$ws.Range("Z25").Value = 7.3
$ws.Range("AA25").Value = 11
$ws.Range("AC25").Value = 13
$ws.Range("AD25").Value = 10.5
$ws.Range("AE25").Value = 25
$ws.Range("AG25").Value = 1000
$ws.Range("AH25").Value = 28
$ws.Range("AN25").Value = 3
$ws.Range("AO25").Value = 5.3
$ws.Range("AQ25").Value = 13.5
$ws.Range("AS25").Value = 250
$ws.Range("AT25").Value = 3.15
$ws.Range("AU25").Value = 9
